$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new values
$ws.Range("A2").Value = "Airén"
$ws.Range("B2").Value = 1300

$ws.Range("A3").Value = "New cometa"
$ws.Range("B3").Value = 1300

$ws.Range("A4").Value = "Cometa"
$ws.Range("B4").Value = 1300

# Remove row 5 entirely so the used range shrinks to A1:B4
$ws.Rows.Item(5).Delete()
